# Auto-generated edit script applying numeric updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 131.54546
$ws.Range("I33").Value = 105.875
$ws.Range("K33").Value = 105.875
$ws.Range("M33").Value = 123.125

$ws.Range("H135").Value = 934.7222
$ws.Range("I135").Value = 814
$ws.Range("J135").Value = 1248.6
$ws.Range("K135").Value = 7326
$ws.Range("L135").Value = 11237.4
$ws.Range("M135").Value = -4791
$ws.Range("N135").Value = -16307.4

$ws.Range("H137").Value = 5393.778
$ws.Range("I137").Value = 1648.3334
$ws.Range("K137").Value = 4945.0002
$ws.Range("M137").Value = -2395.0002

$ws.Range("H138").Value = 2021.0416
$ws.Range("I138").Value = 710.44446
$ws.Range("J138").Value = 5952.8335
$ws.Range("K138").Value = 2131.33338
$ws.Range("L138").Value = 17858.5005
$ws.Range("M138").Value = 3008.66662
$ws.Range("N138").Value = -28138.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 22856.334
$ws.Range("I38").Value = 3529.5
$ws.Range("J38").Value = 61510
$ws.Range("K38").Value = 3529.5
$ws.Range("L38").Value = 61510
$ws.Range("M38").Value = -3062.5
$ws.Range("N38").Value = -62444

$ws.Range("H61").Value = 1536.1538
$ws.Range("J61").Value = 2233
$ws.Range("L61").Value = 2233
$ws.Range("N61").Value = -2657

$ws.Range("H74").Value = 2719
$ws.Range("I74").Value = 2697.2273
$ws.Range("J74").Value = 2838.75
$ws.Range("K74").Value = 2697.2273
$ws.Range("L74").Value = 2838.75
$ws.Range("M74").Value = -1823.2273
$ws.Range("N74").Value = -4586.75

$ws.Range("H77").Value = 2719
$ws.Range("I77").Value = 2697.2273
$ws.Range("J77").Value = 2838.75
$ws.Range("K77").Value = 13486.1365
$ws.Range("L77").Value = 14193.75
$ws.Range("M77").Value = -9118.136500000001
$ws.Range("N77").Value = -22929.75

$ws.Range("H132").Value = 2160.8635
$ws.Range("I132").Value = 1924
$ws.Range("J132").Value = 2792.5
$ws.Range("K132").Value = 5772
$ws.Range("L132").Value = 8377.5
$ws.Range("M132").Value = -3242
$ws.Range("N132").Value = -13437.5

$ws.Range("H136").Value = 1536.1538
$ws.Range("J136").Value = 2233
$ws.Range("L136").Value = 6699
$ws.Range("N136").Value = -11799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3669.24
$ws.Range("I134").Value = 1310.0435
$ws.Range("K134").Value = 3930.1305
$ws.Range("M134").Value = -1395.1305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 95.04761999999999
$ws.Range("I7").Value = 39.666668
$ws.Range("K7").Value = 39.666668
$ws.Range("M7").Value = 73.333332

$ws.Range("H16").Value = 966.5
$ws.Range("I16").Value = 1349.5
$ws.Range("J16").Value = 583.5
$ws.Range("K16").Value = 1349.5
$ws.Range("L16").Value = 583.5
$ws.Range("M16").Value = -1062.5
$ws.Range("N16").Value = -1157.5

$ws.Range("H58").Value = 3479.0833
$ws.Range("I58").Value = 2375.9
$ws.Range("K58").Value = 2375.9
$ws.Range("M58").Value = -2172.9

$ws.Range("H59").Value = 41691.25
$ws.Range("I59").Value = 20000
$ws.Range("K59").Value = 20000
$ws.Range("M59").Value = -18855

$ws.Range("H113").Value = 966.5
$ws.Range("I113").Value = 1349.5
$ws.Range("J113").Value = 583.5
$ws.Range("K113").Value = 1349.5
$ws.Range("L113").Value = 583.5
$ws.Range("M113").Value = 820.5
$ws.Range("N113").Value = -4923.5

$ws.Range("H122").Value = 650
$ws.Range("I122").Value = 650
$ws.Range("K122").Value = 1950
$ws.Range("M122").Value = 500

$ws.Range("H125").Value = 80766.664
$ws.Range("J125").Value = 80766.664
$ws.Range("L125").Value = 80766.664
$ws.Range("N125").Value = -85686.664

$ws.Range("H134").Value = 2260
$ws.Range("I134").Value = 1396
$ws.Range("J134").Value = 6580
$ws.Range("K134").Value = 4188
$ws.Range("L134").Value = 19740
$ws.Range("M134").Value = -1653
$ws.Range("N134").Value = -24810

$ws.Range("H136").Value = 3479.0833
$ws.Range("I136").Value = 2375.9
$ws.Range("K136").Value = 7127.700000000001
$ws.Range("M136").Value = -4577.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1854999.9
$ws.Range("I4").Value = 1854999.9
$ws.Range("K4").Value = 5564999.699999999
$ws.Range("M4").Value = -5564887.699999999

$ws.Range("H5").Value = 1691.9333
$ws.Range("I5").Value = 1373.091
$ws.Range("J5").Value = 2568.75
$ws.Range("K5").Value = 4119.272999999999
$ws.Range("L5").Value = 7706.25
$ws.Range("M5").Value = -4007.272999999999
$ws.Range("N5").Value = -7930.25

$ws.Range("H122").Value = 486.125
$ws.Range("I122").Value = 355.8125
$ws.Range("J122").Value = 746.75
$ws.Range("K122").Value = 3202.3125
$ws.Range("L122").Value = 6720.75
$ws.Range("M122").Value = -752.3125
$ws.Range("N122").Value = -11620.75

$ws.Range("H135").Value = 1691.9333
$ws.Range("I135").Value = 1373.091
$ws.Range("J135").Value = 2568.75
$ws.Range("K135").Value = 12357.819
$ws.Range("L135").Value = 23118.75
$ws.Range("M135").Value = -9822.819
$ws.Range("N135").Value = -28188.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 38905.965
$ws.Range("I132").Value = 51961
$ws.Range("J132").Value = 4636.5
$ws.Range("K132").Value = 155883
$ws.Range("L132").Value = 13909.5
$ws.Range("M132").Value = -153353
$ws.Range("N132").Value = -18969.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2103
$ws.Range("I22").Value = 351
$ws.Range("J22").Value = 2541
$ws.Range("K22").Value = 351
$ws.Range("L22").Value = 2541
$ws.Range("M22").Value = -56
$ws.Range("N22").Value = -3131

$ws.Range("H27").Value = 2103
$ws.Range("I27").Value = 351
$ws.Range("J27").Value = 2541
$ws.Range("K27").Value = 351
$ws.Range("L27").Value = 2541
$ws.Range("M27").Value = -244
$ws.Range("N27").Value = -2755

$ws.Range("H46").Value = 1718.8889
$ws.Range("I46").Value = 518.75
$ws.Range("K46").Value = 518.75
$ws.Range("M46").Value = -330.75

$ws.Range("H55").Value = 803.0625
$ws.Range("I55").Value = 896.2143
$ws.Range("K55").Value = 896.2143
$ws.Range("M55").Value = -723.2143

$ws.Range("H68").Value = 3821.6365
$ws.Range("I68").Value = 2925
$ws.Range("J68").Value = 4897.6
$ws.Range("K68").Value = 2925
$ws.Range("L68").Value = 4897.6
$ws.Range("M68").Value = -2176
$ws.Range("N68").Value = -6395.6

$ws.Range("H71").Value = 3821.6365
$ws.Range("I71").Value = 2925
$ws.Range("J71").Value = 4897.6
$ws.Range("K71").Value = 14625
$ws.Range("L71").Value = 24488
$ws.Range("M71").Value = -10881
$ws.Range("N71").Value = -31976

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 51000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 51000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 51000
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = -51298

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

$ws.Range("H136").Value = 2347.275
$ws.Range("I136").Value = 1735.1936
$ws.Range("K136").Value = 5205.5808
$ws.Range("M136").Value = -2655.5808
